# Generate Report for Handback
# Updates the handback-status report with a freshly-generated run:
#   - file 633aae42-7a24-4a74-bfb8-dcec0e9408be.md -> 2079575f-9c64-4cf3-9b54-60ce6558a8de.md
#   - file a1c5d4d9-4745-493f-b88e-b67490793fcb.md -> ffff298396d1-b8ac-4717-8ca6-8ecd68b28ece.md
#   - refreshed xlf hashes/timestamps for the new handoff/handback cycle

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.md"
$ws.Range("B2").Value = "e2e\2079575f-9c64-4cf3-9b54-60ce6558a8de.md"
$ws.Range("G2").Value = "2016-08-27 21:01:11"

$ws.Range("A3").Value = "ffff298396d1-b8ac-4717-8ca6-8ecd68b28ece.md"
$ws.Range("B3").Value = "e2e\ffff298396d1-b8ac-4717-8ca6-8ecd68b28ece.md"
$ws.Range("G3").Value = "2016-08-27 21:01:11"

# ---------------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.md"
$ws.Range("I2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.md"
$ws.Range("G2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-27 21:01:06"
$ws.Range("J2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-27 21:01:30"

$ws.Range("A3").Value = "ffff298396d1-b8ac-4717-8ca6-8ecd68b28ece.md"
$ws.Range("I3").Value = "ffff298396d1-b8ac-4717-8ca6-8ecd68b28ece.md"
$ws.Range("G3").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-27 21:01:06"
$ws.Range("J3").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-27 21:01:30"

# ---------------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.md"
$ws.Range("I2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.md"
$ws.Range("G2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.de-de.xlf"
$ws.Range("H2").Value = "2016-08-27 21:01:11"
$ws.Range("J2").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.de-de.xlf"
$ws.Range("K2").Value = "2016-08-27 21:01:36"

$ws.Range("A3").Value = "ffff298396d1-b8ac-4717-8ca6-8ecd68b28ece.md"
$ws.Range("I3").Value = "ffff298396d1-b8ac-4717-8ca6-8ecd68b28ece.md"
$ws.Range("G3").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.de-de.xlf"
$ws.Range("H3").Value = "2016-08-27 21:01:11"
$ws.Range("J3").Value = "2079575f-9c64-4cf3-9b54-60ce6558a8de.533eabdaaf1f0b3c4a5a5cdd59e00d45aa56b4f3.de-de.xlf"
$ws.Range("K3").Value = "2016-08-27 21:01:36"
